# Fix misaligned Michigan orgs: the "tags" sheet was missing an
# "African-American" key/value pair that belongs right after "PAC"
# (row 4) and before "Climate" (row 5). Every row from the old
# "Climate" row onward was shifted up by one, which misaligned all
# the tag rows below it (including the Michigan row). Insert the
# missing row to restore the correct alignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tags")

# Insert a new row above the current row 5 ("Climate"), shifting
# it (and everything below) down by one.
$ws.Range("A5").EntireRow.Insert()

# Populate the newly inserted row with the missing tag.
$ws.Range("A5").Value = "African-American"
$ws.Range("B5").Value = "African-American"
